# Apply the data edits described by the diff to the active workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Su (t/sq.m.) column (M) and "Su from Ncor" column (N) for rows 18-31 ---

# M18:M31 all become 0
$ws.Range("M18:M31").Value = 0

# N18, N19, N20 change from 26/35/47 to 50; N21:N31 are already 50.
$ws.Range("N18:N31").Value = 50

# --- Update the view state (scroll position and selection) for the sheet ---
$ws.Range("R27").Select()

$window = $excel.ActiveWindow
$window.ScrollRow = 9
$window.ScrollColumn = 1
